$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Normalisation")

# Row 7 - B7: bestemming_id -> destination_id
$ws.Range("B7").Value = "destination_id"

# Row 13 - remove D13 (border_value_variable), add F13 (border_value)
$ws.Range("D13").Value = ""
$ws.Range("F13").Value = "border_value"

# Row 14 - remove D14 (border_value_number)
$ws.Range("D14").Value = ""

# Row 10 - add F10: BORDER_PER_PRODUCT
$ws.Range("F10").Value = "BORDER_PER_PRODUCT"
$ws.Range("F10").Font.Bold = $true

# Row 11 - add F11 (product_id) and G11 (*=> product)
$ws.Range("F11").Value = "product_id"
$ws.Range("G11").Value = "*=> product"
$ws.Range("G11").Font.Bold = $true

# Row 12 - add F12 (variable_id) and G12 (*=> variable)
$ws.Range("F12").Value = "variable_id"
$ws.Range("G12").Value = "*=> variable"
$ws.Range("G12").Font.Bold = $true

# Row 17 - H17: STABILISATIES_PER_CARGO -> STABILISATIONS_PER_CARGO
$ws.Range("H17").Value = "STABILISATIONS_PER_CARGO"

# Row 17 - F17: exceedings_PER_CARGO -> EXCEEDINGS_PER_CARGO
$ws.Range("F17").Value = "EXCEEDINGS_PER_CARGO"

# Update selection to F22 to match the recorded active cell in the diff
$ws.Range("F22").Select()
